{"js": "// \"minor change in report\" -- apply the three visible edits made to the\n// report text (the remaining hunks in the source diff are proofErr\n// spell/grammar-check markup and Word re-save artifacts that carry no\n// visible/textual change, so they are intentionally not reproduced here).\n\nconst body = context.document.body;\n\n// 1) \"...for that particular year the UFO was sighted\" -> drop \"particular \"\nconst particularResults = body.search(\"particular year\", { matchCase: false, matchWholeWord: false });\nawait context.sync();\nif (particularResults.items.length > 0) {\n  particularResults.items[0].insertText(\"year\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Typo fix: \"...based on the year n which it happened...\" -> \"year in which\"\nconst typoResults = body.search(\"year n which\", { matchCase: false, matchWholeWord: false });\nawait context.sync();\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"year in which\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Bold the file names \"meteorites.json\" and \"meteor.py\"\nconst meteoritesJsonResults = body.search(\"meteorites.json\", { matchCase: true, matchWholeWord: false });\nconst meteorPyResults = body.search(\"meteor.py\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (meteoritesJsonResults.items.length > 0) {\n  meteoritesJsonResults.items[0].font.bold = true;\n}\nif (meteorPyResults.items.length > 0) {\n  meteorPyResults.items[0].font.bold = true;\n}\nawait context.sync();\n\n// 4) The \"_GoBack\" bookmark (Word's \"last edit location\" marker) follows the\n//    edit -- it moves from the old \"...number of UFO[_GoBack] sightings...\"\n//    spot to wrap the newly bolded \"meteorites.json\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\nconst bookmarkTarget = body.search(\"meteorites.json\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (bookmarkTarget.items.length > 0) {\n  bookmarkTarget.items[0].insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"minor change in report\" -- apply the three visible edits made to the\n# report text (the remaining hunks in the source diff are proofErr\n# spell/grammar-check markup and Word re-save artifacts that carry no\n# visible/textual change, so they are intentionally not reproduced here).\n\n$d = $word.ActiveDocument\n\n# 1) \"...for that particular year the UFO was sighted\" -> drop \"particular \"\n$r1 = $d.Content\n$r1.Find.Text = \"particular \"\nif ($r1.Find.Execute()) {\n    $r1.Text = \"\"\n}\n\n# 2) Typo fix: \"...based on the year n which it happened...\" -> \"year in which\"\n$r2 = $d.Content\n$r2.Find.Text = \"year n which\"\n$r2.Find.Replacement.Text = \"year in which\"\n$r2.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 3) Bold the file names \"meteorites.json\" and \"meteor.py\"\n$r3 = $d.Content\n$r3.Find.Text = \"meteorites.json\"\nif ($r3.Find.Execute()) {\n    $r3.Bold = 1\n}\n\n$r4 = $d.Content\n$r4.Find.Text = \"meteor.py\"\nif ($r4.Find.Execute()) {\n    $r4.Bold = 1\n}\n\n# 4) The \"_GoBack\" bookmark (Word's \"last edit location\" marker) follows the\n#    edit -- it moves from the old \"...number of UFO[_GoBack] sightings...\"\n#    spot to wrap the newly bolded \"meteorites.json\".\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$r5 = $d.Content\n$r5.Find.Text = \"meteorites.json\"\nif ($r5.Find.Execute()) {\n    $d.Bookmarks.Add(\"_GoBack\", $r5)\n}\n"}
